# issue #5: add legislator_id, name, date into dataframe
#
# The source workbook's "股票" (stocks) sheet gained three trailing
# columns: date / legislator_name / legislator_id. Re-create that by
# copying the existing header/data formatting onto the new cells (so no
# spurious new styles get created) and writing the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)   # 股票

$legislatorName = "李昆澤"
$legislatorId = 1327
$reportDate = "2012-04-27"

# ---- Header row (row 1): H1=date, I1=legislator_name, J1=legislator_id
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)
$ws.Cells.Item(1, 8).Value = "date"

$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)
$ws.Cells.Item(1, 9).Value = "legislator_name"

$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 10).PasteSpecial(-4122)
$ws.Cells.Item(1, 10).Value = "legislator_id"

# ---- Data row (row 2): H2=date value, I2=legislator name, J2=legislator id
# H2 holds an ISO-looking date string that must stay plain text, not get
# auto-coerced into a date serial: force text format first.
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = $reportDate
# Re-copy the plain data-row formatting over the cell so it matches its
# siblings (clears the temporary "@" text format picked up above without
# disturbing the text we just stored).
$ws.Cells.Item(2, 7).Copy()
$ws.Cells.Item(2, 8).PasteSpecial(-4122)

$ws.Cells.Item(2, 7).Copy()
$ws.Cells.Item(2, 9).PasteSpecial(-4122)
$ws.Cells.Item(2, 9).Value = $legislatorName

$ws.Cells.Item(2, 7).Copy()
$ws.Cells.Item(2, 10).PasteSpecial(-4122)
$ws.Cells.Item(2, 10).Value = $legislatorId
